$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: report title / date range (advance one week) ---
# "Volume 30   Number  17" -> "Volume 30   Number  18"
$ws.Range("A8").Value = "Volume 30   Number  18"
# "Report Covering the Week  4/24/2023  Through  4/30/2023"
#   -> "Report Covering the Week  5/1/2023  Through  5/7/2023"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Crime-complaint table body (rows 14-30) ---
# Cells that become the textual placeholders "0" / "***.*" are populated by
# copying (value + style) from existing reference cells C14 ("0", style 14)
# and E14 ("***.*", style 14) so the shared-string + style indices line up
# exactly with the rest of the sheet.
# Cells that become numeric (where they used to hold the text placeholder)
# get an explicit NumberFormat so they pick up the same style as their
# numeric neighbours (#,##0 for counts, #,##0.0;"-"#,##0.0 for percentages).

    # Row 14
    $ws.Range("L14").Value = 0
    # Row 15
    $ws.Range("D15").Value = 1
    $ws.Range("D15").NumberFormat = '#,##0'
    $ws.Range("E15").Value = -100
    $ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("F15").Value = 2
    $ws.Range("G15").Value = 2
    $ws.Range("H15").Value = 0
    $ws.Range("J15").Value = 7
    $ws.Range("K15").Value = 0
    $ws.Range("M15").Value = 40
    # Row 16
    $ws.Range("C16").Value = 4
    $ws.Range("D16").Value = 4
    $ws.Range("E16").Value = 0
    $ws.Range("F16").Value = 13
    $ws.Range("H16").Value = 0
    $ws.Range("I16").Value = 68
    $ws.Range("J16").Value = 54
    $ws.Range("K16").Value = 25.925925925925
    $ws.Range("L16").Value = 83.783783783783
    $ws.Range("M16").Value = -20
    $ws.Range("N16").Value = -80.232558139534
    # Row 17
    $ws.Range("C17").Value = 6
    $ws.Range("D17").Value = 6
    $ws.Range("E17").Value = 0
    $ws.Range("F17").Value = 26
    $ws.Range("G17").Value = 21
    $ws.Range("H17").Value = 23.809523809523
    $ws.Range("I17").Value = 95
    $ws.Range("J17").Value = 83
    $ws.Range("K17").Value = 14.457831325301
    $ws.Range("L17").Value = 14.457831325301
    $ws.Range("M17").Value = 23.376623376623
    $ws.Range("N17").Value = 4.395604395604
    # Row 18
    $ws.Range("C18").Value = 4
    $ws.Range("D18").Value = 6
    $ws.Range("E18").Value = -33.333333333333
    $ws.Range("F18").Value = 12
    $ws.Range("G18").Value = 16
    $ws.Range("H18").Value = -25
    $ws.Range("I18").Value = 65
    $ws.Range("J18").Value = 82
    $ws.Range("K18").Value = -20.731707317073
    $ws.Range("L18").Value = -9.722222222222
    $ws.Range("M18").Value = -53.900709219858
    $ws.Range("N18").Value = -90.767045454545
    # Row 19
    $ws.Range("C19").Value = 19
    $ws.Range("D19").Value = 8
    $ws.Range("E19").Value = 137.5
    $ws.Range("F19").Value = 51
    $ws.Range("G19").Value = 49
    $ws.Range("H19").Value = 4.081632653061
    $ws.Range("I19").Value = 217
    $ws.Range("J19").Value = 216
    $ws.Range("K19").Value = 0.462962962962
    $ws.Range("L19").Value = 32.317073170731
    $ws.Range("M19").Value = 68.217054263565
    $ws.Range("N19").Value = 5.339805825242
    # Row 20
    $ws.Range("D20").Value = 13
    $ws.Range("E20").Value = -30.769230769230
    $ws.Range("F20").Value = 34
    $ws.Range("G20").Value = 23
    $ws.Range("H20").Value = 47.826086956521
    $ws.Range("I20").Value = 117
    $ws.Range("J20").Value = 87
    $ws.Range("K20").Value = 34.482758620689
    $ws.Range("L20").Value = 112.727272727273
    $ws.Range("M20").Value = -4.098360655737
    $ws.Range("N20").Value = -91.183119819140
    # Row 21
    $ws.Range("C21").Value = 42
    $ws.Range("D21").Value = 38
    $ws.Range("E21").Value = 10.526315789473
    $ws.Range("F21").Value = 138
    $ws.Range("G21").Value = 124
    $ws.Range("H21").Value = 11.290322580645
    $ws.Range("I21").Value = 571
    $ws.Range("J21").Value = 529
    $ws.Range("K21").Value = 7.939508506616
    $ws.Range("L21").Value = 36.276849642004
    $ws.Range("M21").Value = 1.964285714285
    $ws.Range("N21").Value = -78.773234200743
    # Row 22
    $ws.Range("C14").Copy($ws.Range("D22"))
    $ws.Range("E14").Copy($ws.Range("E22"))
    $ws.Range("F22").Value = 2
    $ws.Range("H22").Value = 0
    # Row 24
    $ws.Range("C24").Value = 18
    $ws.Range("D24").Value = 23
    $ws.Range("E24").Value = -21.739130434782
    $ws.Range("F24").Value = 69
    $ws.Range("G24").Value = 113
    $ws.Range("H24").Value = -38.938053097345
    $ws.Range("I24").Value = 426
    $ws.Range("J24").Value = 456
    $ws.Range("K24").Value = -6.578947368421
    $ws.Range("L24").Value = -10.878661087866
    $ws.Range("M24").Value = 13.903743315508
    # Row 25
    $ws.Range("C25").Value = 14
    $ws.Range("D25").Value = 12
    $ws.Range("E25").Value = 16.666666666666
    $ws.Range("F25").Value = 50
    $ws.Range("G25").Value = 50
    $ws.Range("I25").Value = 184
    $ws.Range("J25").Value = 168
    $ws.Range("K25").Value = 9.523809523809
    $ws.Range("L25").Value = 30.496453900709
    $ws.Range("M25").Value = -28.957528957529
    # Row 26
    $ws.Range("C14").Copy($ws.Range("C26"))
    $ws.Range("D26").Value = 2
    $ws.Range("D26").NumberFormat = '#,##0'
    $ws.Range("E26").Value = -100
    $ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("F26").Value = 3
    $ws.Range("G26").Value = 4
    $ws.Range("H26").Value = -25
    $ws.Range("J26").Value = 13
    $ws.Range("K26").Value = -23.076923076923
    # Row 27
    $ws.Range("C14").Copy($ws.Range("C27"))
    $ws.Range("C14").Copy($ws.Range("D27"))
    $ws.Range("E14").Copy($ws.Range("E27"))
    $ws.Range("F27").Value = 6
    $ws.Range("G27").Value = 5
    $ws.Range("H27").Value = 20
    $ws.Range("L27").Value = 50
    # Row 28
    $ws.Range("N28").Value = -75
    # Row 29
    $ws.Range("N29").Value = -71.428571428571
    # Row 30
    $ws.Range("C14").Copy($ws.Range("D30"))
    $ws.Range("E14").Copy($ws.Range("E30"))

